# Generate Report for Handback
#
# The localization handback pipeline has produced (and handed back) the
# de-de and zh-cn packages. Update the status report:
#   - flip the "Status" cell everywhere from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - record the generated target (.md) file and handback (.xlf) file for
#     each language, with a hyperlink on the target-file cell
#   - record the handback timestamp for each language
#   - widen the columns that now hold the longer values

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# "Status" is shown on the Overview sheet (once per language column) and on
# each language sheet's own Status cell - update every occurrence so the
# underlying text changes everywhere it is displayed.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

$mdFileName = "f9293860-4eeb-46f8-823e-dd5eb0e8d0b1.md"
$mdHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f22ff513ff21769dce9ed1619036ea69119b30a0/e2e/f9293860-4eeb-46f8-823e-dd5eb0e8d0b1.md"

# ---- zh-cn: record the newly generated target + handback files ----
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Range("J2").Value = "f9293860-4eeb-46f8-823e-dd5eb0e8d0b1.23b205b032d6c51cb84835bbb80ede0934c7ae69.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 07:04:52"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    $mdHyperlinkUrl,
    "",
    "",
    $mdFileName
)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276

# ---- de-de: record the newly generated target + handback files ----
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Range("J2").Value = "f9293860-4eeb-46f8-823e-dd5eb0e8d0b1.23b205b032d6c51cb84835bbb80ede0934c7ae69.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 07:04:59"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    $mdHyperlinkUrl,
    "",
    "",
    $mdFileName
)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276

# ---- Widen columns to fit the longer values now populated above ----
# (ColumnWidth is quantized to whole pixels by the host, same as real
# Excel; use the column-center width so it lands on the intended pixel.)
$wsOverview.Range("E:F").ColumnWidth = 175 / 6
$wsZhCn.Range("C:C").ColumnWidth = 175 / 6
$wsDeDe.Range("C:C").ColumnWidth = 175 / 6

$wsZhCn.Range("I:I").ColumnWidth = 235 / 6
$wsZhCn.Range("J:J").ColumnWidth = 235 / 6
$wsDeDe.Range("I:I").ColumnWidth = 235 / 6
$wsDeDe.Range("J:J").ColumnWidth = 235 / 6
